# Script: separate "All" category into "Other" on main.model sheet,
# and remove the "fitness" row from the domestic-grazer block on the
# regressions sheet (begun separating invertebrates from vertebrates).

$wb = $excel.ActiveWorkbook

# --- main.model sheet: rename "All" rows to "Other" ---
$wsModel = $wb.Worksheets.Item("main.model")
$wsModel.Range("A3").Value = "Other"
$wsModel.Range("A14").Value = "Other"

# --- regressions sheet: delete the "fitness" row (row 4) ---
$wsReg = $wb.Worksheets.Item("regressions")
$wsReg.Rows.Item(4).Delete()

# --- restore view/selection state ---
$wsModel.Activate()
$excel.ActiveWindow.ScrollRow = 1
$wsModel.Range("A9").Select()

$wsReg.Activate()
$wsReg.Range("C15").Select()

$wsModel.Activate()
